$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Temperatura (C) and Oxigeno (D) readings for rows 4-15.
$ws.Cells.Item(4, 3).Value  = 34
$ws.Cells.Item(4, 4).Value  = 98

$ws.Cells.Item(5, 3).Value  = 33
$ws.Cells.Item(5, 4).Value  = 99

$ws.Cells.Item(6, 3).Value  = 32
$ws.Cells.Item(6, 4).Value  = 100

$ws.Cells.Item(7, 3).Value  = 31
$ws.Cells.Item(7, 4).Value  = 101

$ws.Cells.Item(8, 3).Value  = 30
$ws.Cells.Item(8, 4).Value  = 102

$ws.Cells.Item(9, 3).Value  = 29
$ws.Cells.Item(9, 4).Value  = 103

$ws.Cells.Item(10, 3).Value = 28
$ws.Cells.Item(10, 4).Value = 104

$ws.Cells.Item(11, 3).Value = 27
$ws.Cells.Item(11, 4).Value = 105

$ws.Cells.Item(12, 3).Value = 26
$ws.Cells.Item(12, 4).Value = 106

$ws.Cells.Item(13, 3).Value = 25
$ws.Cells.Item(13, 4).Value = 107

$ws.Cells.Item(14, 3).Value = 24
$ws.Cells.Item(14, 4).Value = 108

$ws.Cells.Item(15, 3).Value = 23
$ws.Cells.Item(15, 4).Value = 109

# Update the active selection to match the edited range.
[void]$ws.Range("D2:D15").Select()
